$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 597.25
$ws.Range("I6").Value = 678
$ws.Range("J6").Value = 355
$ws.Range("K6").Value = 2034
$ws.Range("L6").Value = 1065
$ws.Range("M6").Value = -1922
$ws.Range("N6").Value = -1289

$ws.Range("H9").Value = 113
$ws.Range("I9").Value = 70
$ws.Range("J9").Value = 199
$ws.Range("K9").Value = 70
$ws.Range("L9").Value = 199
$ws.Range("M9").Value = 99
$ws.Range("N9").Value = -537

$ws.Range("H12").Value = 976.6667
$ws.Range("I12").Value = 1387.25
$ws.Range("J12").Value = 155.5
$ws.Range("K12").Value = 1387.25
$ws.Range("L12").Value = 155.5
$ws.Range("M12").Value = -1217.25
$ws.Range("N12").Value = -495.5

$ws.Range("H38").Value = 99
$ws.Range("I38").Value = 99
$ws.Range("K38").Value = 297
$ws.Range("M38").Value = 75

$ws.Range("H98").Value = 25473.572
$ws.Range("I98").Value = 21933.334
$ws.Range("J98").Value = 28128.75
$ws.Range("K98").Value = 21933.334
$ws.Range("L98").Value = 28128.75
$ws.Range("M98").Value = -20435.334
$ws.Range("N98").Value = -31124.75

$ws.Range("H122").Value = 25473.572
$ws.Range("I122").Value = 21933.334
$ws.Range("J122").Value = 28128.75
$ws.Range("K122").Value = 65800.00199999999
$ws.Range("L122").Value = 84386.25
$ws.Range("M122").Value = -63350.00199999999
$ws.Range("N122").Value = -89286.25

$ws.Range("H141").Value = 3283.4285
$ws.Range("I141").Value = 3283.4285
$ws.Range("K141").Value = 9850.2855
$ws.Range("M141").Value = -4670.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 110
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -332

$ws.Range("H45").Value = 797.2222
$ws.Range("I45").Value = 772
$ws.Range("K45").Value = 772
$ws.Range("M45").Value = -395

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2955.077
$ws.Range("I99").Value = 2834.6667
$ws.Range("K99").Value = 2834.6667
$ws.Range("M99").Value = -1336.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12125.083
$ws.Range("I31").Value = 6416.1665
$ws.Range("J31").Value = 17834
$ws.Range("K31").Value = 6416.1665
$ws.Range("L31").Value = 17834
$ws.Range("M31").Value = -6121.1665
$ws.Range("N31").Value = -18424

$ws.Range("H34").Value = 12125.083
$ws.Range("I34").Value = 6416.1665
$ws.Range("J34").Value = 17834
$ws.Range("K34").Value = 6416.1665
$ws.Range("L34").Value = 17834
$ws.Range("M34").Value = -6214.1665
$ws.Range("N34").Value = -18238

$ws.Range("H51").Value = 29380.2
$ws.Range("J51").Value = 25475.25
$ws.Range("L51").Value = 25475.25
$ws.Range("N51").Value = -26947.25

$ws.Range("H58").Value = 2217
$ws.Range("J58").Value = 2492.5
$ws.Range("L58").Value = 2492.5
$ws.Range("N58").Value = -2898.5

$ws.Range("H60").Value = 24286.666
$ws.Range("I60").Value = 10296
$ws.Range("J60").Value = 31282
$ws.Range("K60").Value = 10296
$ws.Range("L60").Value = 31282
$ws.Range("M60").Value = -9785
$ws.Range("N60").Value = -32304

$ws.Range("H61").Value = 29380.2
$ws.Range("J61").Value = 25475.25
$ws.Range("L61").Value = 25475.25
$ws.Range("N61").Value = -26171.25

$ws.Range("H74").Value = 89044
$ws.Range("J74").Value = 89044
$ws.Range("L74").Value = 89044
$ws.Range("N74").Value = -90792

$ws.Range("H77").Value = 89044
$ws.Range("J77").Value = 89044
$ws.Range("L77").Value = 267132
$ws.Range("N77").Value = -275868

$ws.Range("H86").Value = 5949.6
$ws.Range("I86").Value = 7149.5
$ws.Range("J86").Value = 5149.6665
$ws.Range("K86").Value = 7149.5
$ws.Range("L86").Value = 5149.6665
$ws.Range("M86").Value = -6026.5
$ws.Range("N86").Value = -7395.6665

$ws.Range("H89").Value = 5949.6
$ws.Range("I89").Value = 7149.5
$ws.Range("J89").Value = 5149.6665
$ws.Range("K89").Value = 35747.5
$ws.Range("L89").Value = 25748.3325
$ws.Range("M89").Value = -30131.5
$ws.Range("N89").Value = -36980.3325

$ws.Range("H134").Value = 979.6667
$ws.Range("I134").Value = 924.4286
$ws.Range("K134").Value = 2773.2858
$ws.Range("M134").Value = -238.2857999999997

$ws.Range("H136").Value = 2217
$ws.Range("J136").Value = 2492.5
$ws.Range("L136").Value = 7477.5
$ws.Range("N136").Value = -12577.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 510
$ws.Range("I33").Value = 650.25
$ws.Range("J33").Value = 416.5
$ws.Range("K33").Value = 3901.5
$ws.Range("L33").Value = 2499
$ws.Range("M33").Value = -3618.5
$ws.Range("N33").Value = -3065

$ws.Range("H86").Value = 2900
$ws.Range("I86").Value = 2700
$ws.Range("K86").Value = 8100
$ws.Range("M86").Value = -6914

$ws.Range("H89").Value = 2900
$ws.Range("I89").Value = 2700
$ws.Range("K89").Value = 24300
$ws.Range("M89").Value = -18372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 3140
$ws.Range("I39").Value = 2710
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 2710
$ws.Range("L39").Value = 4000
$ws.Range("M39").Value = -2250
$ws.Range("N39").Value = -4920

$ws.Range("H46").Value = 224933
$ws.Range("I46").Value = 500750
$ws.Range("K46").Value = 500750
$ws.Range("M46").Value = -500562

$ws.Range("H55").Value = 1499.5
$ws.Range("J55").Value = 1499.5
$ws.Range("L55").Value = 1499.5
$ws.Range("N55").Value = -1845.5

$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251

$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 60000
$ws.Range("I52").Value = 60000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 60000
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("M52").Value = -59774

$ws.Range("H96").Value = 22953
$ws.Range("J96").Value = 50250
$ws.Range("L96").Value = 50250
$ws.Range("N96").Value = -52996
